$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,7).Value = 3.336603
$ws.Cells.Item(2,8).Value = 10.009809
$ws.Cells.Item(2,9).Value = 0.2530037693731349
$ws.Cells.Item(2,10).Value = 0.253003769373135
$ws.Cells.Item(2,13).Value = 60.538204
$ws.Cells.Item(2,14).Value = 181.614612
$ws.Cells.Item(2,15).Value = 0.6123615450168176
$ws.Cells.Item(2,16).Value = 0.6123615450168176
$ws.Cells.Item(2,17).Value = 201.991953081012
$ws.Cells.Item(2,18).Value = 1817.927577729108
$ws.Cells.Item(2,19).Value = 0.1549297791084115
$ws.Cells.Item(2,20).Value = 0.1549297791084115

# Row 3
$ws.Cells.Item(3,7).Value = 3.336603
$ws.Cells.Item(3,8).Value = 10.009809
$ws.Cells.Item(3,9).Value = 0.2530037693731349
$ws.Cells.Item(3,10).Value = 0.253003769373135
$ws.Cells.Item(3,15).Value = 0.1096681363892149
$ws.Cells.Item(3,16).Value = 0.1096681363892149
$ws.Cells.Item(3,17).Value = 36.174840239852
$ws.Cells.Item(3,18).Value = 325.573562158668
$ws.Cells.Item(3,19).Value = 0.02774645188659842
$ws.Cells.Item(3,20).Value = 0.02774645188659843

# Row 4
$ws.Cells.Item(4,7).Value = 3.336603
$ws.Cells.Item(4,8).Value = 10.009809
$ws.Cells.Item(4,9).Value = 0.2530037693731349
$ws.Cells.Item(4,10).Value = 0.253003769373135
$ws.Cells.Item(4,13).Value = 8.850437666666666
$ws.Cells.Item(4,14).Value = 26.551313
$ws.Cells.Item(4,15).Value = 0.08952475173586316
$ws.Cells.Item(4,16).Value = 0.08952475173586316
$ws.Cells.Item(4,17).Value = 29.530396869913
$ws.Cells.Item(4,18).Value = 265.773571829217
$ws.Cells.Item(4,19).Value = 0.02265009964136748
$ws.Cells.Item(4,20).Value = 0.02265009964136749

# Row 5
$ws.Cells.Item(5,7).Value = 3.336603
$ws.Cells.Item(5,8).Value = 10.009809
$ws.Cells.Item(5,9).Value = 0.2530037693731349
$ws.Cells.Item(5,10).Value = 0.253003769373135
$ws.Cells.Item(5,13).Value = 1.757142
$ws.Cells.Item(5,14).Value = 5.271426
$ws.Cells.Item(5,15).Value = 0.01777400251143792
$ws.Cells.Item(5,16).Value = 0.01777400251143792
$ws.Cells.Item(5,17).Value = 5.862885268626
$ws.Cells.Item(5,18).Value = 52.765967417634
$ws.Cells.Item(5,19).Value = 0.00449688963224136
$ws.Cells.Item(5,20).Value = 0.00449688963224136

# Row 6
$ws.Cells.Item(6,7).Value = 3.336603
$ws.Cells.Item(6,8).Value = 10.009809
$ws.Cells.Item(6,9).Value = 0.2530037693731349
$ws.Cells.Item(6,10).Value = 0.253003769373135
$ws.Cells.Item(6,13).Value = 16.87263033333333
$ws.Cells.Item(6,14).Value = 50.617891
$ws.Cells.Item(6,15).Value = 0.1706715643466665
$ws.Cells.Item(6,16).Value = 0.1706715643466665
$ws.Cells.Item(6,17).Value = 56.29726898809101
$ws.Cells.Item(6,18).Value = 506.675420892819
$ws.Cells.Item(6,19).Value = 0.04318054910451616
$ws.Cells.Item(6,20).Value = 0.04318054910451617

# Row 7
$ws.Cells.Item(7,9).Value = 0.4389650376240971
$ws.Cells.Item(7,10).Value = 0.4389650376240971
$ws.Cells.Item(7,13).Value = 60.538204
$ws.Cells.Item(7,14).Value = 181.614612
$ws.Cells.Item(7,15).Value = 0.6123615450168176
$ws.Cells.Item(7,16).Value = 0.6123615450168176
$ws.Cells.Item(7,17).Value = 350.4588311220093
$ws.Cells.Item(7,18).Value = 3154.129480098084
$ws.Cells.Item(7,19).Value = 0.2688053086478576
$ws.Cells.Item(7,20).Value = 0.2688053086478576

# Row 8
$ws.Cells.Item(8,9).Value = 0.4389650376240971
$ws.Cells.Item(8,10).Value = 0.4389650376240971
$ws.Cells.Item(8,15).Value = 0.1096681363892149
$ws.Cells.Item(8,16).Value = 0.1096681363892149
$ws.Cells.Item(8,19).Value = 0.04814047761625631
$ws.Cells.Item(8,20).Value = 0.04814047761625632

# Row 9
$ws.Cells.Item(9,9).Value = 0.4389650376240971
$ws.Cells.Item(9,10).Value = 0.4389650376240971
$ws.Cells.Item(9,13).Value = 8.850437666666666
$ws.Cells.Item(9,14).Value = 26.551313
$ws.Cells.Item(9,15).Value = 0.08952475173586316
$ws.Cells.Item(9,16).Value = 0.08952475173586316
$ws.Cells.Item(9,17).Value = 51.23564682523789
$ws.Cells.Item(9,18).Value = 461.120821427141
$ws.Cells.Item(9,19).Value = 0.03929823601402112
$ws.Cells.Item(9,20).Value = 0.03929823601402112

# Row 10
$ws.Cells.Item(10,9).Value = 0.4389650376240971
$ws.Cells.Item(10,10).Value = 0.4389650376240971
$ws.Cells.Item(10,13).Value = 1.757142
$ws.Cells.Item(10,14).Value = 5.271426
$ws.Cells.Item(10,15).Value = 0.01777400251143792
$ws.Cells.Item(10,16).Value = 0.01777400251143792
$ws.Cells.Item(10,17).Value = 10.172186995098
$ws.Cells.Item(10,18).Value = 91.54968295588199
$ws.Cells.Item(10,19).Value = 0.007802165681164141
$ws.Cells.Item(10,20).Value = 0.007802165681164141

# Row 11
$ws.Cells.Item(11,9).Value = 0.4389650376240971
$ws.Cells.Item(11,10).Value = 0.4389650376240971
$ws.Cells.Item(11,13).Value = 16.87263033333333
$ws.Cells.Item(11,14).Value = 50.617891
$ws.Cells.Item(11,15).Value = 0.1706715643466665
$ws.Cells.Item(11,16).Value = 0.1706715643466665
$ws.Cells.Item(11,17).Value = 97.67654000065411
$ws.Cells.Item(11,18).Value = 879.088860005887
$ws.Cells.Item(11,19).Value = 0.07491884966479795
$ws.Cells.Item(11,20).Value = 0.07491884966479795

# Row 12
$ws.Cells.Item(12,7).Value = 1.029432
$ws.Cells.Item(12,8).Value = 3.088296
$ws.Cells.Item(12,9).Value = 0.07805848532574147
$ws.Cells.Item(12,10).Value = 0.07805848532574149
$ws.Cells.Item(12,13).Value = 60.538204
$ws.Cells.Item(12,14).Value = 181.614612
$ws.Cells.Item(12,15).Value = 0.6123615450168176
$ws.Cells.Item(12,16).Value = 0.6123615450168176
$ws.Cells.Item(12,17).Value = 62.319964420128
$ws.Cells.Item(12,18).Value = 560.8796797811519
$ws.Cells.Item(12,19).Value = 0.04780001467574364
$ws.Cells.Item(12,20).Value = 0.04780001467574364

# Row 13
$ws.Cells.Item(13,7).Value = 1.029432
$ws.Cells.Item(13,8).Value = 3.088296
$ws.Cells.Item(13,9).Value = 0.07805848532574147
$ws.Cells.Item(13,10).Value = 0.07805848532574149
$ws.Cells.Item(13,15).Value = 0.1096681363892149
$ws.Cells.Item(13,16).Value = 0.1096681363892149
$ws.Cells.Item(13,17).Value = 11.160913701088
$ws.Cells.Item(13,18).Value = 100.448223309792
$ws.Cells.Item(13,19).Value = 0.008560528615038943
$ws.Cells.Item(13,20).Value = 0.008560528615038946

# Row 14
$ws.Cells.Item(14,7).Value = 1.029432
$ws.Cells.Item(14,8).Value = 3.088296
$ws.Cells.Item(14,9).Value = 0.07805848532574147
$ws.Cells.Item(14,10).Value = 0.07805848532574149
$ws.Cells.Item(14,13).Value = 8.850437666666666
$ws.Cells.Item(14,14).Value = 26.551313
$ws.Cells.Item(14,15).Value = 0.08952475173586316
$ws.Cells.Item(14,16).Value = 0.08952475173586316
$ws.Cells.Item(14,17).Value = 9.110923748071999
$ws.Cells.Item(14,18).Value = 81.998313732648
$ws.Cells.Item(14,19).Value = 0.006988166519664523
$ws.Cells.Item(14,20).Value = 0.006988166519664525

# Row 15
$ws.Cells.Item(15,7).Value = 1.029432
$ws.Cells.Item(15,8).Value = 3.088296
$ws.Cells.Item(15,9).Value = 0.07805848532574147
$ws.Cells.Item(15,10).Value = 0.07805848532574149
$ws.Cells.Item(15,13).Value = 1.757142
$ws.Cells.Item(15,14).Value = 5.271426
$ws.Cells.Item(15,15).Value = 0.01777400251143792
$ws.Cells.Item(15,16).Value = 0.01777400251143792
$ws.Cells.Item(15,17).Value = 1.808858203344
$ws.Cells.Item(15,18).Value = 16.279723830096
$ws.Cells.Item(15,19).Value = 0.001387411714218769
$ws.Cells.Item(15,20).Value = 0.001387411714218769

# Row 16
$ws.Cells.Item(16,7).Value = 1.029432
$ws.Cells.Item(16,8).Value = 3.088296
$ws.Cells.Item(16,9).Value = 0.07805848532574147
$ws.Cells.Item(16,10).Value = 0.07805848532574149
$ws.Cells.Item(16,13).Value = 16.87263033333333
$ws.Cells.Item(16,14).Value = 50.617891
$ws.Cells.Item(16,15).Value = 0.1706715643466665
$ws.Cells.Item(16,16).Value = 0.1706715643466665
$ws.Cells.Item(16,17).Value = 17.369225589304
$ws.Cells.Item(16,18).Value = 156.323030303736
$ws.Cells.Item(16,19).Value = 0.01332236380107561
$ws.Cells.Item(16,20).Value = 0.01332236380107561

# Row 17
$ws.Cells.Item(17,7).Value = 1.675087666666667
$ws.Cells.Item(17,8).Value = 5.025263
$ws.Cells.Item(17,9).Value = 0.1270164576658104
$ws.Cells.Item(17,10).Value = 0.1270164576658104
$ws.Cells.Item(17,13).Value = 60.538204
$ws.Cells.Item(17,14).Value = 181.614612
$ws.Cells.Item(17,15).Value = 0.6123615450168176
$ws.Cells.Item(17,16).Value = 0.6123615450168176
$ws.Cells.Item(17,17).Value = 101.4067988825507
$ws.Cells.Item(17,18).Value = 912.661189942956
$ws.Cells.Item(17,19).Value = 0.07777999425879886
$ws.Cells.Item(17,20).Value = 0.07777999425879886

# Row 18
$ws.Cells.Item(18,7).Value = 1.675087666666667
$ws.Cells.Item(18,8).Value = 5.025263
$ws.Cells.Item(18,9).Value = 0.1270164576658104
$ws.Cells.Item(18,10).Value = 0.1270164576658104
$ws.Cells.Item(18,15).Value = 0.1096681363892149
$ws.Cells.Item(18,16).Value = 0.1096681363892149
$ws.Cells.Item(18,17).Value = 18.16099449931955
$ws.Cells.Item(18,18).Value = 163.448950493876
$ws.Cells.Item(18,19).Value = 0.01392965820296903
$ws.Cells.Item(18,20).Value = 0.01392965820296903

# Row 19
$ws.Cells.Item(19,7).Value = 1.675087666666667
$ws.Cells.Item(19,8).Value = 5.025263
$ws.Cells.Item(19,9).Value = 0.1270164576658104
$ws.Cells.Item(19,10).Value = 0.1270164576658104
$ws.Cells.Item(19,13).Value = 8.850437666666666
$ws.Cells.Item(19,14).Value = 26.551313
$ws.Cells.Item(19,15).Value = 0.08952475173586316
$ws.Cells.Item(19,16).Value = 0.08952475173586316
$ws.Cells.Item(19,17).Value = 14.82525898003544
$ws.Cells.Item(19,18).Value = 133.427330820319
$ws.Cells.Item(19,19).Value = 0.01137111683890045
$ws.Cells.Item(19,20).Value = 0.01137111683890045

# Row 20
$ws.Cells.Item(20,7).Value = 1.675087666666667
$ws.Cells.Item(20,8).Value = 5.025263
$ws.Cells.Item(20,9).Value = 0.1270164576658104
$ws.Cells.Item(20,10).Value = 0.1270164576658104
$ws.Cells.Item(20,13).Value = 1.757142
$ws.Cells.Item(20,14).Value = 5.271426
$ws.Cells.Item(20,15).Value = 0.01777400251143792
$ws.Cells.Item(20,16).Value = 0.01777400251143792
$ws.Cells.Item(20,17).Value = 2.943366892782
$ws.Cells.Item(20,18).Value = 26.490302035038
$ws.Cells.Item(20,19).Value = 0.002257590837546062
$ws.Cells.Item(20,20).Value = 0.002257590837546062

# Row 21
$ws.Cells.Item(21,7).Value = 1.675087666666667
$ws.Cells.Item(21,8).Value = 5.025263
$ws.Cells.Item(21,9).Value = 0.1270164576658104
$ws.Cells.Item(21,10).Value = 0.1270164576658104
$ws.Cells.Item(21,13).Value = 16.87263033333333
$ws.Cells.Item(21,14).Value = 50.617891
$ws.Cells.Item(21,15).Value = 0.1706715643466665
$ws.Cells.Item(21,16).Value = 0.1706715643466665
$ws.Cells.Item(21,17).Value = 28.26313497559255
$ws.Cells.Item(21,18).Value = 254.368214780333
$ws.Cells.Item(21,19).Value = 0.021678097527596
$ws.Cells.Item(21,20).Value = 0.021678097527596

# Row 22
$ws.Cells.Item(22,5).Value = 3
$ws.Cells.Item(22,6).Value = 1
$ws.Cells.Item(22,7).Value = 1.357782666666667
$ws.Cells.Item(22,8).Value = 4.073348
$ws.Cells.Item(22,9).Value = 0.102956250011216
$ws.Cells.Item(22,10).Value = 0.102956250011216
$ws.Cells.Item(22,13).Value = 60.538204
$ws.Cells.Item(22,14).Value = 181.614612
$ws.Cells.Item(22,15).Value = 0.6123615450168176
$ws.Cells.Item(22,16).Value = 0.6123615450168176
$ws.Cells.Item(22,17).Value = 82.19772406233068
$ws.Cells.Item(22,18).Value = 739.779516560976
$ws.Cells.Item(22,19).Value = 0.063046448326006
$ws.Cells.Item(22,20).Value = 0.063046448326006

# Row 23
$ws.Cells.Item(23,5).Value = 3
$ws.Cells.Item(23,6).Value = 1
$ws.Cells.Item(23,7).Value = 1.357782666666667
$ws.Cells.Item(23,8).Value = 4.073348
$ws.Cells.Item(23,9).Value = 0.102956250011216
$ws.Cells.Item(23,10).Value = 0.102956250011216
$ws.Cells.Item(23,15).Value = 0.1096681363892149
$ws.Cells.Item(23,16).Value = 0.1096681363892149
$ws.Cells.Item(23,17).Value = 14.72083165036622
$ws.Cells.Item(23,18).Value = 132.487484853296
$ws.Cells.Item(23,19).Value = 0.01129102006835215
$ws.Cells.Item(23,20).Value = 0.01129102006835215

# Row 24
$ws.Cells.Item(24,5).Value = 3
$ws.Cells.Item(24,6).Value = 1
$ws.Cells.Item(24,7).Value = 1.357782666666667
$ws.Cells.Item(24,8).Value = 4.073348
$ws.Cells.Item(24,9).Value = 0.102956250011216
$ws.Cells.Item(24,10).Value = 0.102956250011216
$ws.Cells.Item(24,13).Value = 8.850437666666666
$ws.Cells.Item(24,14).Value = 26.551313
$ws.Cells.Item(24,15).Value = 0.08952475173586316
$ws.Cells.Item(24,16).Value = 0.08952475173586316
$ws.Cells.Item(24,17).Value = 12.01697085621378
$ws.Cells.Item(24,18).Value = 108.152737705924
$ws.Cells.Item(24,19).Value = 0.009217132721909575
$ws.Cells.Item(24,20).Value = 0.009217132721909575

# Row 25
$ws.Cells.Item(25,5).Value = 3
$ws.Cells.Item(25,6).Value = 1
$ws.Cells.Item(25,7).Value = 1.357782666666667
$ws.Cells.Item(25,8).Value = 4.073348
$ws.Cells.Item(25,9).Value = 0.102956250011216
$ws.Cells.Item(25,10).Value = 0.102956250011216
$ws.Cells.Item(25,13).Value = 1.757142
$ws.Cells.Item(25,14).Value = 5.271426
$ws.Cells.Item(25,15).Value = 0.01777400251143792
$ws.Cells.Item(25,16).Value = 0.01777400251143792
$ws.Cells.Item(25,17).Value = 2.385816950472
$ws.Cells.Item(25,18).Value = 21.472352554248
$ws.Cells.Item(25,19).Value = 0.001829944646267584
$ws.Cells.Item(25,20).Value = 0.001829944646267584

# Row 26
$ws.Cells.Item(26,5).Value = 3
$ws.Cells.Item(26,6).Value = 1
$ws.Cells.Item(26,7).Value = 1.357782666666667
$ws.Cells.Item(26,8).Value = 4.073348
$ws.Cells.Item(26,9).Value = 0.102956250011216
$ws.Cells.Item(26,10).Value = 0.102956250011216
$ws.Cells.Item(26,13).Value = 16.87263033333333
$ws.Cells.Item(26,14).Value = 50.617891
$ws.Cells.Item(26,15).Value = 0.1706715643466665
$ws.Cells.Item(26,16).Value = 0.1706715643466665
$ws.Cells.Item(26,17).Value = 22.90936500767423
$ws.Cells.Item(26,18).Value = 206.184285069068
$ws.Cells.Item(26,19).Value = 0.01757170424868074
$ws.Cells.Item(26,20).Value = 0.01757170424868074

